$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.986.11"
$ws.Range("E2").Value = "  -1.33%  "
$ws.Range("D3").Value = "2.932.97"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'376.82"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("D6").Value = "'102.06"
$ws.Range("E6").Value = "  -4.11%  "
$ws.Range("D7").Value = "'0.534"
$ws.Range("E7").Value = "  -2.70%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "'0.581"
$ws.Range("E9").Value = "  -3.65%  "
$ws.Range("D10").Value = "'36.43"
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "'0.139"
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -1.52%  "
$ws.Range("D13").Value = "3.390.16"
$ws.Range("E13").Value = "  -2.68%  "
$ws.Range("D14").Value = "'17.89"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "2.930.40"
$ws.Range("E16").Value = "  -2.88%  "
$ws.Range("E17").Value = "  -0.93%  "
$ws.Range("D18").Value = "50.908.99"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("D19").Value = "'3.15"
$ws.Range("E19").Value = "  -8.53%  "
$ws.Range("D20").Value = "'7.11"
$ws.Range("E20").Value = "  -4.83%  "
$ws.Range("D21").Value = "'12.44"
$ws.Range("E21").Value = "  -5.04%  "
$ws.Range("E22").Value = "  -2.02%  "
$ws.Range("D23").Value = "'68.03"
$ws.Range("E23").Value = "  -1.31%  "
$ws.Range("D24").Value = "'261.19"
$ws.Range("E24").Value = "  -1.13%  "
$ws.Range("E25").Value = "  +1.69%  "
$ws.Range("D26").Value = "'8.19"
$ws.Range("E26").Value = "  +8.98%  "
$ws.Range("D27").Value = "'7.51"
$ws.Range("E27").Value = "  +2.25%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  -3.61%  "
$ws.Range("E30").Value = "  +5.26%  "
$ws.Range("D31").Value = "'25.50"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("D32").Value = "'9.77"
$ws.Range("E32").Value = "  -2.20%  "
$ws.Range("D33").Value = "'50.58"
$ws.Range("E33").Value = "  -1.18%  "
$ws.Range("D34").Value = "'33.74"
$ws.Range("E34").Value = "  -3.63%  "
$ws.Range("D35").Value = "'0.0452"
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("E36").Value = "  -3.01%  "
$ws.Range("E37").Value = "  -0.14%  "
$ws.Range("E38").Value = "  -4.96%  "
$ws.Range("D39").Value = "'2.55"
$ws.Range("E39").Value = "  -3.61%  "
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").Value = "'16.24"
$ws.Range("E41").Value = "  -7.88%  "
$ws.Range("E42").Value = "  -5.33%  "
$ws.Range("D43").Value = "'120.95"
$ws.Range("E43").Value = "  -2.93%  "
$ws.Range("D44").Value = "'21.06"
$ws.Range("E44").Value = "  -5.97%  "
$ws.Range("D45").Value = "'2.05"
$ws.Range("E45").Value = "  -1.73%  "
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "2.001.59"
$ws.Range("E48").Value = "  -2.81%  "
$ws.Range("D49").Value = "'3.20"
$ws.Range("E49").Value = "  -2.79%  "
$ws.Range("D50").Value = "'0.0341"
$ws.Range("E50").Value = "  -3.98%  "
$ws.Range("E51").Value = "  -4.47%  "
